$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column I, mirroring column B/E data with some different values
$ws.Range("I2").Value = 3
$ws.Range("I3").Value = "P1"
$ws.Range("I4").Value = 10
$ws.Range("I5").Value = "P2"
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = "P3"
$ws.Range("I8").Value = 10
$ws.Range("I9").Value = 1
$ws.Range("I10").Value = 2
$ws.Range("I11").Value = 3
$ws.Range("I12").Value = 5

# Add new rows 13 and 14 in column E
$ws.Range("E13").Value = 3
$ws.Range("E14").Value = 5

# Update selection to I2
$ws.Range("I2").Select()
